# Etisalat Mashari Al Afasi - add RBT content start/expire date columns (StartDate/ExpireDate)
# for rows 2-22, plus header styling (bordered/centered/bold) and date-format cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New column widths for E (StartDate) and F (ExpireDate)
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 12.86
$ws.Columns.Item(6).ColumnWidth = 17.1667

# ---------------------------------------------------------------------------
# 2) Header cells E1 / F1
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "StartDate"
$ws.Range("F1").Value = "ExpireDate"

$headerRng = $ws.Range("E1:F1")
foreach ($addr in @("E1", "F1")) {
    $hc = $ws.Range($addr)
    $hc.Borders(7).LineStyle = 1
    $hc.Borders(10).LineStyle = 1
}
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3) Data rows 2-22: StartDate = 2020-10-02 (44106), ExpireDate = 2020-10-27 (44131)
# ---------------------------------------------------------------------------
$startSerial = 44106
$expireSerial = 44131

for ($r = 2; $r -le 22; $r++) {
    $ws.Range("E$r").Value = $startSerial
    $ws.Range("F$r").Value = $expireSerial
}

$dateRng = $ws.Range("E2:F22")
$dateRng.NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# 4) Sheet view: drop the saved scroll position / old selection, select F22
# ---------------------------------------------------------------------------
$ws.Range("F22").Select()

Write-Host "RBT StartDate/ExpireDate columns added"
